$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").ClearContents()
$ws.Range("A8").Select()
